# Added analysis of prices
# Append newly-observed price-history rows to each product sheet.
# Dates/times are written as literal text (not auto-converted to Excel
# date/time serials), matching the existing inlineStr-text convention
# used throughout the workbook.

$wb = $excel.ActiveWorkbook

function Add-PriceRow {
    param($ws, $row, $date, $time, $price)

    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $date
    $dateCell.Style = "Normal"

    $timeCell = $ws.Cells.Item($row, 2)
    $timeCell.NumberFormat = "@"
    $timeCell.Value = $time
    $timeCell.Style = "Normal"

    $ws.Cells.Item($row, 3).Value = $price
}

# Sheet 1: Zotac Gaming RTX 3060  (A1:C3 -> A1:C4)
$ws1 = $wb.Worksheets.Item("Zotac Gaming RTX 3060")
Add-PriceRow $ws1 4 "2025-05-02" "13:07:33" 25999

# Sheet 2: Infinix Note 50x 5G+  (A1:C3 -> A1:C5)
$ws2 = $wb.Worksheets.Item("Infinix Note 50x 5G+")
Add-PriceRow $ws2 4 "2025-05-02" "13:07:39" 11499
Add-PriceRow $ws2 5 "2025-05-18" "22:27:42" 11499

# Sheet 3: DJI MIC 2  (A1:C3 -> A1:C4)
$ws3 = $wb.Worksheets.Item("DJI MIC 2")
Add-PriceRow $ws3 4 "2025-05-02" "13:07:47" 20900

# Sheet 4: Nothing Phone (3a) - 256 GB  (A1:C2 -> A1:C4)
$ws4 = $wb.Worksheets.Item("Nothing Phone (3a) - 256 GB")
Add-PriceRow $ws4 3 "2025-05-02" "13:07:53" 26999
Add-PriceRow $ws4 4 "2025-05-18" "22:27:53" 26999
